# Vega Monumental Concepción - Berenjena: weekly fruit/veg price update.
# A new observation (week) is inserted as row 19 and every existing row
# from the old row 19 down to the old row 42 shifts down by one (to rows
# 20-43), preserving its data unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 19 - this pushes rows 19:42 down
# to 20:43 and grows the used range to A1:R43.
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly observation.
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = "Vega Monumental Concepción"
$ws.Range("C19").Value = "Bíobío"
$ws.Range("D19").Value = 44495
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 100112001
$ws.Range("G19").Value = "Berenjena"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 8500
$ws.Range("M19").Value = 8222
$ws.Range("N19").Value = "$/caja 60 unidades"
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 137
$ws.Range("Q19").Value = 60
$ws.Range("R19").Value = "Hortaliza"
